$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 311; this shifts the existing rows 311-343
# down to 312-344 and updates the sheet dimension automatically.
$ws.Rows("311").Insert()

# Populate the newly inserted row 311 with the new record.
$ws.Range("A311").Value = 4
$ws.Range("B311").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C311").Value = "Los Lagos"
$ws.Range("D311").Value = 45142
$ws.Range("E311").Value = 10
$ws.Range("F311").Value = 100112009
$ws.Range("G311").Value = "Acelga"
$ws.Range("H311").Value = "Sin especificar"
$ws.Range("I311").Value = "Primera"
$ws.Range("J311").Value = 75
$ws.Range("K311").Value = 10000
$ws.Range("L311").Value = 10000
$ws.Range("M311").Value = 10000
$ws.Range("N311").Value = "$/docena de atados (12 kilos)"
$ws.Range("O311").Value = "Región de La Araucanía"
$ws.Range("P311").Value = 833
$ws.Range("Q311").Value = 12
$ws.Range("R311").Value = "Hortaliza"
